$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1599.8334
$ws.Range("I51").Value = 1999.5
$ws.Range("J51").Value = 1400
$ws.Range("K51").Value = 1999.5
$ws.Range("L51").Value = 1400
$ws.Range("M51").Value = -1515.5
$ws.Range("N51").Value = -2368

$ws.Range("H98").Value = 1513.238
$ws.Range("I98").Value = 1076.9231
$ws.Range("J98").Value = 2222.25
$ws.Range("K98").Value = 1076.9231
$ws.Range("L98").Value = 2222.25
$ws.Range("M98").Value = 421.0769
$ws.Range("N98").Value = -5218.25

$ws.Range("H111").Value = 142858850
$ws.Range("I111").Value = 142858850
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 428576550
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -428573483

$ws.Range("H113").Value = 1772.2727
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1849.5
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1849.5
$ws.Range("M113").Value = 2254
$ws.Range("N113").Value = -8357.5

$ws.Range("H122").Value = 1513.238
$ws.Range("I122").Value = 1076.9231
$ws.Range("J122").Value = 2222.25
$ws.Range("K122").Value = 3230.7693
$ws.Range("L122").Value = 6666.75
$ws.Range("M122").Value = -780.7692999999999
$ws.Range("N122").Value = -11566.75

$ws.Range("H131").Value = 26645.375
$ws.Range("I131").Value = 32130.188
$ws.Range("J131").Value = 4706.125
$ws.Range("K131").Value = 96390.564
$ws.Range("L131").Value = 14118.375
$ws.Range("M131").Value = -91350.564
$ws.Range("N131").Value = -24198.375

$ws.Range("H135").Value = 473.125
$ws.Range("I135").Value = 405.17392
$ws.Range("J135").Value = 2036
$ws.Range("K135").Value = 3646.56528
$ws.Range("L135").Value = 18324
$ws.Range("M135").Value = -1111.56528
$ws.Range("N135").Value = -23394

$ws.Range("H137").Value = 1710.8889
$ws.Range("I137").Value = 1348.0741
$ws.Range("J137").Value = 2799.3333
$ws.Range("K137").Value = 4044.2223
$ws.Range("L137").Value = 8397.999899999999
$ws.Range("M137").Value = -1494.2223
$ws.Range("N137").Value = -13497.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3303.1428
$ws.Range("I2").Value = 3303.1428
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3303.1428
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -3190.1428
$ws.Range("N2").Value = ""

$ws.Range("H74").Value = 4959.407
$ws.Range("I74").Value = 5073.231
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 5073.231
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -4199.231
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 4959.407
$ws.Range("I77").Value = 5073.231
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 25366.155
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -20998.155
$ws.Range("N77").Value = -18736

$ws.Range("H110").Value = 4400
$ws.Range("I110").Value = 1800
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 1800
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 245
$ws.Range("N110").Value = -11090

$ws.Range("H116").Value = 3303.1428
$ws.Range("I116").Value = 3303.1428
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3303.1428
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1009.1428
$ws.Range("N116").Value = ""

$ws.Range("H132").Value = 2321.5386
$ws.Range("I132").Value = 1797.2632
$ws.Range("J132").Value = 3744.5715
$ws.Range("K132").Value = 5391.7896
$ws.Range("L132").Value = 11233.7145
$ws.Range("M132").Value = -2861.7896
$ws.Range("N132").Value = -16293.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3303.1428
$ws.Range("I3").Value = 3303.1428
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3303.1428
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3189.1428
$ws.Range("N3").Value = ""

$ws.Range("H107").Value = 2789.9
$ws.Range("I107").Value = 2789.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2789.9
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -869.9000000000001
$ws.Range("N107").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 949.5
$ws.Range("I107").Value = 778.1667
$ws.Range("J107").Value = 1206.5
$ws.Range("K107").Value = 778.1667
$ws.Range("L107").Value = 1206.5
$ws.Range("M107").Value = 1141.8333
$ws.Range("N107").Value = -5046.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 707.6
$ws.Range("I5").Value = 547.8
$ws.Range("J5").Value = 867.4
$ws.Range("K5").Value = 1643.4
$ws.Range("L5").Value = 2602.2
$ws.Range("M5").Value = -1531.4
$ws.Range("N5").Value = -2826.2

$ws.Range("H23").Value = 230.14285
$ws.Range("I23").Value = 274.5
$ws.Range("J23").Value = 222.75
$ws.Range("K23").Value = 823.5
$ws.Range("L23").Value = 668.25
$ws.Range("M23").Value = -588.5
$ws.Range("N23").Value = -1138.25

$ws.Range("H86").Value = 1034.5
$ws.Range("I86").Value = 802
$ws.Range("J86").Value = 1081
$ws.Range("K86").Value = 2406
$ws.Range("L86").Value = 3243
$ws.Range("M86").Value = -1220
$ws.Range("N86").Value = -5615

$ws.Range("H89").Value = 1034.5
$ws.Range("I89").Value = 802
$ws.Range("J89").Value = 1081
$ws.Range("K89").Value = 7218
$ws.Range("L89").Value = 9729
$ws.Range("M89").Value = -1290
$ws.Range("N89").Value = -21585

$ws.Range("H98").Value = 221
$ws.Range("I98").Value = 237.8
$ws.Range("J98").Value = 207
$ws.Range("K98").Value = 713.4000000000001
$ws.Range("L98").Value = 621
$ws.Range("M98").Value = 784.5999999999999
$ws.Range("N98").Value = -3617

$ws.Range("H135").Value = 707.6
$ws.Range("I135").Value = 547.8
$ws.Range("J135").Value = 867.4
$ws.Range("K135").Value = 4930.2
$ws.Range("L135").Value = 7806.599999999999
$ws.Range("M135").Value = -2395.2
$ws.Range("N135").Value = -12876.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2219.3914
$ws.Range("I102").Value = 2149.7896
$ws.Range("J102").Value = 2550
$ws.Range("K102").Value = 2149.7896
$ws.Range("L102").Value = 2550
$ws.Range("M102").Value = -527.7896000000001
$ws.Range("N102").Value = -5794

$ws.Range("H107").Value = 182.84616
$ws.Range("I107").Value = 204
$ws.Range("J107").Value = 66.5
$ws.Range("K107").Value = 204
$ws.Range("L107").Value = 66.5
$ws.Range("M107").Value = 1716
$ws.Range("N107").Value = -3906.5

$ws.Range("H113").Value = 13659.889
$ws.Range("I113").Value = 2227.8
$ws.Range("J113").Value = 27950
$ws.Range("K113").Value = 2227.8
$ws.Range("L113").Value = 27950
$ws.Range("M113").Value = -57.80000000000018
$ws.Range("N113").Value = -32290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7595.7856
$ws.Range("I61").Value = 10328.6
$ws.Range("J61").Value = 763.75
$ws.Range("K61").Value = 10328.6
$ws.Range("L61").Value = 763.75
$ws.Range("M61").Value = -10126.6
$ws.Range("N61").Value = -1167.75

$ws.Range("H113").Value = 7595.7856
$ws.Range("I113").Value = 10328.6
$ws.Range("J113").Value = 763.75
$ws.Range("K113").Value = 10328.6
$ws.Range("L113").Value = 763.75
$ws.Range("M113").Value = -8158.6
$ws.Range("N113").Value = -5103.75

$ws.Range("H122").Value = 1899.7142
$ws.Range("I122").Value = 1124.75
$ws.Range("J122").Value = 2933
$ws.Range("K122").Value = 3374.25
$ws.Range("L122").Value = 8799
$ws.Range("M122").Value = -924.25
$ws.Range("N122").Value = -13699

$ws.Range("H132").Value = 6793.2354
$ws.Range("I132").Value = 5581.8335
$ws.Range("J132").Value = 9700.6
$ws.Range("K132").Value = 16745.5005
$ws.Range("L132").Value = 29101.8
$ws.Range("M132").Value = -14215.5005
$ws.Range("N132").Value = -34161.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 562.25
$ws.Range("I107").Value = 562.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1686.75
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 233.25

$ws.Range("H113").Value = 339.81818
$ws.Range("I113").Value = 334.1111
$ws.Range("J113").Value = 365.5
$ws.Range("K113").Value = 1002.3333
$ws.Range("L113").Value = 1096.5
$ws.Range("M113").Value = 1167.6667
$ws.Range("N113").Value = -5436.5

$ws.Range("H126").Value = 819.7727
$ws.Range("I126").Value = 891
$ws.Range("J126").Value = 577.6
$ws.Range("K126").Value = 2673
$ws.Range("L126").Value = 1732.8
$ws.Range("M126").Value = -203
$ws.Range("N126").Value = -6672.8
